# Upload dataset statistics: refresh the bug-count samples on the
# "multiple bugs" sheets and restore each sheet's last-used cell
# selection, matching the author's re-export of the workbook.

$wb = $excel.ActiveWorkbook

$wsSingleSystem = $wb.Worksheets.Item("system_top_1_single_bug")
$wsMultiSystem  = $wb.Worksheets.Item("system_top_1_multiple_bugs")
$wsSingleTopN   = $wb.Worksheets.Item("top_n_single_bug")
$wsMultiTopN    = $wb.Worksheets.Item("top_n_multiple_bugs")

# --- Updated sample data -------------------------------------------------

# system_top_1_multiple_bugs!B4:G4  (ZipMe, ExamDB counts change)
$wsMultiSystem.Range("B4").Value = 45
$wsMultiSystem.Range("E4").Value = 26

# top_n_multiple_bugs!B4:F4  (Top-1..Top-5 counts change)
$wsMultiTopN.Range("B4").Value = 207
$wsMultiTopN.Range("C4").Value = 637
$wsMultiTopN.Range("D4").Value = 1001
$wsMultiTopN.Range("E4").Value = 1096
$wsMultiTopN.Range("F4").Value = 1126

# --- Restore each sheet view's last active-cell selection ----------------

$wsSingleSystem.Range("C34").Select()
$wsMultiSystem.Range("V13").Select()
$wsSingleTopN.Range("G4").Select()
$wsMultiTopN.Range("B10").Select()

# Leave the final selection on the sheet that was active/tabbed in the
# source workbook.
$wsMultiTopN.Activate()
